## Automatische test-sync: 2025-07-23 13:57:50
## Adds a new mail-log row (row 3) to the "Logs" sheet and bumps the
## corresponding count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Logs" sheet: append the new row of data
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Range("D3").Value = "Openingstijden / Locatie"
$logs.Range("E3").Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F3").Value = "2025-07-23 13:57:01"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Ja"
$logs.Range("J3").Value = "Ja"

# Undo Excel's automatic row-height expansion caused by the multi-line
# cell above, so the row keeps the sheet's default height (matches the
# original file, which has no explicit row heights).
$logs.Rows.Item(3).AutoFit()

# Extend the existing conditional-formatting rules (previously only
# covering row 2) so that they also cover the freshly added row 3.
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))
$logs.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J3"))

# ---------------------------------------------------------------
# 2. "Dashboard" sheet: bump the tally for "Openingstijden / Locatie"
# ---------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 2
